$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is brand new; give its date cell (A19) the same format as the other
# date cells in column A (copy format only, from A2 which already uses it).
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(19, 1).PasteSpecial(-4122)

# Target data for rows 2..19 (date serial, y_0, y_0_forecast, y_1, y_1_forecast)
$data = @(
    @(39400, 2007, 11.13090654781819, 2008, 9.396507498425466),
    @(39765, 2008, 4.672550446571067, 2009, 4.422525088127283),
    @(40130, 2009, -14.45332333832743, 2010, -2.928447329610073),
    @(40494, 2010, 8.600536527919633, 2011, 6.303897256856628),
    @(40862, 2011, 10.25770250047622, 2012, 10.22374275635105),
    @(41228, 2012, 4.639893381363169, 2013, 8.174613408931286),
    @(41592, 2013, 0.3058963467304165, 2014, 2.429116709932622),
    @(41957, 2014, 4.068173739091874, 2015, 3.9413000500929),
    @(42321, 2015, 4.984288257750213, 2016, 4.188839638544284),
    @(42689, 2016, 1.878184267712912, 2017, 2.514670279852349),
    @(43053, 2017, 4.695933104194339, 2018, 4.5579527192392),
    @(43418, 2018, 4.892602738886098, 2019, 1.957202207503861),
    @(43783, 2019, 0.8049382522247184, 2020, 2.267257846564918),
    @(44159, 2020, -8.784173899737169, 2021, 2.199380357735481),
    @(44525, 2021, 5.110501195359984, 2022, 0.3515918738370427),
    @(44890, 2022, 5.120680133083599, 2023, 5.934275247805543),
    @(45254, 2023, -0.5532735011319234, 2024, -1.846917864698006),
    @(45618, 2024, -1.069674659641462, 2025, -0.7986414110784379)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row = $row + 1
}
